# Update the "Förändrad" (Changed) date column (C) for all data rows.
# The diff shows every C2:C200 cell value moving from the date serial
# 45204 (2023-10-05) to 45205 (2023-10-06), while everything else
# (formatting, other columns) stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 200 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value2 = 45205
    }
}
